$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "Problems"/"X" marker column from D to E to make room for the new "Url" column F ---
$ws.Range("D1").Cut($ws.Range("E1"))

$ws.Range("D23").Cut($ws.Range("E23"))
$ws.Range("D23").Clear()

$ws.Range("D37:D39").Cut($ws.Range("E37:E39"))

# --- New hyperlink column with SPDR holdings download links ---
# (Added in the same order the original author did: rows 6-39 first, then 5,4,3,2)
$ws.Hyperlinks.Add($ws.Range("F6"), "https://us.spdrs.com/site-content/xls/XLC_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://us.spdrs.com/site-content/xls/XLE_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://us.spdrs.com/site-content/xls/XLF_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://us.spdrs.com/site-content/xls/XLI_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://us.spdrs.com/site-content/xls/XLK_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://us.spdrs.com/site-content/xls/XLP_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://us.spdrs.com/site-content/xls/XLRE_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F13"), "https://us.spdrs.com/site-content/xls/XLU_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F14"), "https://us.spdrs.com/site-content/xls/XLV_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F15"), "https://us.spdrs.com/site-content/xls/XLY_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F16"), "https://us.spdrs.com/site-content/xls/XME_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F17"), "https://us.spdrs.com/site-content/xls/XNTK_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F18"), "https://us.spdrs.com/site-content/xls/XBI_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F19"), "https://us.spdrs.com/site-content/xls/XAR_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F20"), "https://us.spdrs.com/site-content/xls/XOP_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F21"), "https://us.spdrs.com/site-content/xls/XPH_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F22"), "https://us.spdrs.com/site-content/xls/DIA_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F23"), "https://us.spdrs.com/site-content/xls/XWEB_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F24"), "https://us.spdrs.com/site-content/xls/XTL_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F25"), "https://us.spdrs.com/site-content/xls/XSW_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F26"), "https://us.spdrs.com/site-content/xls/XTH_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F27"), "https://us.spdrs.com/site-content/xls/XSD_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F28"), "https://us.spdrs.com/site-content/xls/XRT_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F29"), "https://us.spdrs.com/site-content/xls/XRE_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F30"), "https://us.spdrs.com/site-content/xls/KIE_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F31"), "https://us.spdrs.com/site-content/xls/KCE_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F32"), "https://us.spdrs.com/site-content/xls/KBE_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F33"), "https://us.spdrs.com/site-content/xls/SYE_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F34"), "https://us.spdrs.com/site-content/xls/SPY_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F35"), "https://us.spdrs.com/site-content/xls/SLY_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F36"), "https://us.spdrs.com/site-content/xls/MDY_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F37"), "https://us.spdrs.com/site-content/xls/XTN_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F38"), "https://us.spdrs.com/site-content/xls/XHE_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F39"), "https://us.spdrs.com/site-content/xls/XHS_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://us.spdrs.com/site-content/xls/XLB_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://us.spdrs.com/site-content/xls/XITK_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://us.spdrs.com/site-content/xls/XHB_All_Holdings.xls")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://us.spdrs.com/site-content/xls/XES_All_Holdings.xls")

# --- New header (added last, after the hyperlink column was populated) ---
$ws.Range("F1").Value = "Url"

# --- Match final selection shown in the workbook ---
$ws.Range("F2").Select() | Out-Null
